$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data between row 66 and row 67 (columns F:V) ---
# Row 66 becomes the Chaves vs Gil Vicente match (previously row 67)
# Row 67 becomes the SC Farense vs Vizela match (previously row 66)
$ws.Range("F66").Value = "Chaves"
$ws.Range("G66").Value = 4
$ws.Range("H66").Value = "Gil Vicente"
$ws.Range("I66").Value = 2
$ws.Range("J66").Value = 2.62
$ws.Range("K66").Value = "02/10/2023 20:42"
$ws.Range("L66").Value = 2.74
$ws.Range("M66").Value = "07/10/2023 16:29"
$ws.Range("N66").Value = 3.45
$ws.Range("O66").Value = "02/10/2023 20:42"
$ws.Range("P66").Value = 3.65
$ws.Range("Q66").Value = "07/10/2023 16:25"
$ws.Range("R66").Value = 2.73
$ws.Range("S66").Value = "02/10/2023 20:42"
$ws.Range("T66").Value = 2.58
$ws.Range("U66").Value = "07/10/2023 16:29"
$ws.Range("V66").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/chaves-gil-vicente/K4BKKZh1/"

$ws.Range("F67").Value = "SC Farense"
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = "Vizela"
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 2.29
$ws.Range("K67").Value = "02/10/2023 07:12"
$ws.Range("L67").Value = 2.32
$ws.Range("M67").Value = "07/10/2023 16:02"
$ws.Range("N67").Value = 3.4
$ws.Range("O67").Value = "02/10/2023 07:12"
$ws.Range("P67").Value = 3.47
$ws.Range("Q67").Value = "07/10/2023 15:49"
$ws.Range("R67").Value = 3.31
$ws.Range("S67").Value = "02/10/2023 07:12"
$ws.Range("T67").Value = 3.24
$ws.Range("U67").Value = "07/10/2023 16:03"
$ws.Range("V67").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/sc-farense-vizela/OY1Asc0E/"

# --- Swap match data between row 76 and row 77 (columns F:V) ---
# Row 76 becomes the Vitoria Guimaraes vs Chaves match (previously row 77)
# Row 77 becomes the Benfica vs Casa Pia match (previously row 76)
$ws.Range("F76").Value = "Vitoria Guimaraes"
$ws.Range("G76").Value = 5
$ws.Range("H76").Value = "Chaves"
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 1.62
$ws.Range("K76").Value = "11/10/2023 14:42"
$ws.Range("L76").Value = 1.69
$ws.Range("M76").Value = "28/10/2023 18:58"
$ws.Range("N76").Value = 4.26
$ws.Range("O76").Value = "11/10/2023 14:42"
$ws.Range("P76").Value = 4
$ws.Range("Q76").Value = "28/10/2023 18:58"
$ws.Range("R76").Value = 5.63
$ws.Range("S76").Value = "11/10/2023 14:42"
$ws.Range("T76").Value = 5.29
$ws.Range("U76").Value = "28/10/2023 18:58"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/vitoria-guimaraes-chaves/8vH9wlat/"

$ws.Range("F77").Value = "Benfica"
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = "Casa Pia"
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = 1.22
$ws.Range("K77").Value = "11/10/2023 14:42"
$ws.Range("L77").Value = 1.22
$ws.Range("M77").Value = "28/10/2023 18:55"
$ws.Range("N77").Value = 7.32
$ws.Range("O77").Value = "11/10/2023 14:42"
$ws.Range("P77").Value = 6.95
$ws.Range("Q77").Value = "28/10/2023 18:58"
$ws.Range("R77").Value = 13.22
$ws.Range("S77").Value = "11/10/2023 14:42"
$ws.Range("T77").Value = 14.17
$ws.Range("U77").Value = "28/10/2023 18:58"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/benfica-casa-pia/GWtkzFhl/"

# --- Append new row 101: Chaves vs Vizela (15/11-01/12/2023) ---
# Copy formatting (styles) from the last existing data row (100) to the new row 101
$ws.Range("A100:V100").Copy()
$ws.Range("A101:V101").PasteSpecial(-4122)

$ws.Range("A101").Value = 100
$ws.Range("B101").Value = "portugal"
$ws.Range("C101").Value = "liga-portugal"
$ws.Range("D101").Value = "2023-2024"
$ws.Range("E101").Value = 45261.89583333334
$ws.Range("F101").Value = "Chaves"
$ws.Range("G101").Value = 2
$ws.Range("H101").Value = "Vizela"
$ws.Range("I101").Value = 1
$ws.Range("J101").Value = 2.76
$ws.Range("K101").Value = "15/11/2023 15:12"
$ws.Range("L101").Value = 2.94
$ws.Range("M101").Value = "01/12/2023 21:27"
$ws.Range("N101").Value = 3.42
$ws.Range("O101").Value = "15/11/2023 15:12"
$ws.Range("P101").Value = 3.28
$ws.Range("Q101").Value = "01/12/2023 21:27"
$ws.Range("R101").Value = 2.49
$ws.Range("S101").Value = "15/11/2023 15:12"
$ws.Range("T101").Value = 2.51
$ws.Range("U101").Value = "01/12/2023 20:24"
$ws.Range("V101").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/chaves-vizela/Ichs1rc4/"
